$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A21 with refined timestamp value (tiny precision correction from diff)
$ws.Range("A21").Value = 45876.79187491898

# Add new row 22 data
$ws.Range("A22").Value = 45876.83356961413
$ws.Range("B22").Value = 2025
$ws.Range("C22").Value = 28
$ws.Range("D22").Value = 14.72
$ws.Range("E22").Value = 89.76000000000001
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 3.65
$ws.Range("H22").Value = "ESE"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "20:00:20"

# Copy style from A21 to A22 (date format)
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
